$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date corrections on existing rows (2015/2016 trade dates) ---
# Row 4: TradeDate moved from 2015-05-24 to 2016-05-25 (EffectiveDate recalculates via formula)
$ws.Cells.Item(4, 2).Value = 42515

# Row 5: TradeDate shifted by +1 day (2015-12-19 -> 2015-12-20)
$ws.Cells.Item(5, 2).Value = 42358

# Row 6: TradeDate shifted by +1 day (2016-03-20 -> 2016-03-21)
$ws.Cells.Item(6, 2).Value = 42450

# Row 7: TradeDate shifted by +1 day (2015-10-11 -> 2015-10-12)
$ws.Cells.Item(7, 2).Value = 42289

# --- Rate correction on row 82 (Rate 1.2635 -> 1.5, SecondaryAmount recalculates via formula) ---
$ws.Cells.Item(82, 9).Value = 1.5

# --- New trade rows 102-123 (22 additional blotter entries) ---
$ws.Cells.Item(102, 1).Value = 101
$ws.Cells.Item(102, 2).Value = 42504
$ws.Cells.Item(102, 3).Formula = "=SUM(B102+1)"
$ws.Cells.Item(102, 4).Value = 'Swap'
$ws.Cells.Item(102, 5).Value = 'EUR'
$ws.Cells.Item(102, 6).Value = -50000000
$ws.Cells.Item(102, 7).Value = 'USD'
$ws.Cells.Item(102, 8).Formula = "=SUM(F102*I102*-1)"
$ws.Cells.Item(102, 9).Value = 1.1304
$ws.Cells.Item(102, 9).Style = "Normal"
$ws.Cells.Item(102, 11).Value = 'Nomura'
$ws.Cells.Item(102, 12).Value = 42509
$ws.Cells.Item(102, 13).Value = 'Kelley Babb'
$ws.Cells.Item(102, 14).Value = 'Approved'

$ws.Cells.Item(103, 1).Value = 102
$ws.Cells.Item(103, 2).Value = 42513
$ws.Cells.Item(103, 3).Formula = "=SUM(B103+1)"
$ws.Cells.Item(103, 4).Value = 'Spot'
$ws.Cells.Item(103, 5).Value = 'GBP'
$ws.Cells.Item(103, 6).Value = 1000000
$ws.Cells.Item(103, 7).Value = 'USD'
$ws.Cells.Item(103, 8).Formula = "=SUM(F103*I103*-1)"
$ws.Cells.Item(103, 9).Value = 1.4352
$ws.Cells.Item(103, 9).Style = "Normal"
$ws.Cells.Item(103, 11).Value = 'Lloyds TSB'
$ws.Cells.Item(103, 12).Value = 42509
$ws.Cells.Item(103, 13).Value = 'Sarai Pilgrim'
$ws.Cells.Item(103, 14).Value = 'Approved'

$ws.Cells.Item(104, 1).Value = 103
$ws.Cells.Item(104, 2).Value = 42519
$ws.Cells.Item(104, 3).Formula = "=SUM(B104+1)"
$ws.Cells.Item(104, 4).Value = 'Spot'
$ws.Cells.Item(104, 5).Value = 'GBP'
$ws.Cells.Item(104, 6).Value = 10000
$ws.Cells.Item(104, 7).Value = 'USD'
$ws.Cells.Item(104, 8).Formula = "=SUM(F104*I104*-1)"
$ws.Cells.Item(104, 9).Value = 1.4352
$ws.Cells.Item(104, 9).Style = "Normal"
$ws.Cells.Item(104, 11).Value = 'Goldman Sachs'
$ws.Cells.Item(104, 12).Value = 42509
$ws.Cells.Item(104, 13).Value = 'Stacee Dreiling'
$ws.Cells.Item(104, 14).Value = 'Approved'

$ws.Cells.Item(105, 1).Value = 104
$ws.Cells.Item(105, 2).Value = 42532
$ws.Cells.Item(105, 3).Formula = "=SUM(B105+1)"
$ws.Cells.Item(105, 4).Value = 'Forward'
$ws.Cells.Item(105, 5).Value = 'EUR'
$ws.Cells.Item(105, 6).Value = -10000
$ws.Cells.Item(105, 7).Value = 'JPY'
$ws.Cells.Item(105, 8).Formula = "=SUM(F105*I105*-1)"
$ws.Cells.Item(105, 9).Value = 123.9075
$ws.Cells.Item(105, 9).Style = "Normal"
$ws.Cells.Item(105, 11).Value = 'Morgan Stanley'
$ws.Cells.Item(105, 12).Value = 42509
$ws.Cells.Item(105, 13).Value = 'Cecil Staab'
$ws.Cells.Item(105, 14).Value = 'Approved'

$ws.Cells.Item(106, 1).Value = 105
$ws.Cells.Item(106, 2).Value = 42572
$ws.Cells.Item(106, 3).Formula = "=SUM(B106+1)"
$ws.Cells.Item(106, 4).Value = 'Forward'
$ws.Cells.Item(106, 5).Value = 'EUR'
$ws.Cells.Item(106, 6).Value = -50000000
$ws.Cells.Item(106, 7).Value = 'GBP'
$ws.Cells.Item(106, 8).Formula = "=SUM(F106*I106*-1)"
$ws.Cells.Item(106, 9).Value = 0.7876
$ws.Cells.Item(106, 9).Style = "Normal"
$ws.Cells.Item(106, 11).Value = 'Barcap'
$ws.Cells.Item(106, 12).Value = 42509
$ws.Cells.Item(106, 13).Value = 'Louella Spiker'
$ws.Cells.Item(106, 14).Value = 'Approved'

$ws.Cells.Item(107, 1).Value = 106
$ws.Cells.Item(107, 2).Value = 42586
$ws.Cells.Item(107, 3).Formula = "=SUM(B107+1)"
$ws.Cells.Item(107, 4).Value = 'Swap'
$ws.Cells.Item(107, 5).Value = 'EUR'
$ws.Cells.Item(107, 6).Value = 50000000
$ws.Cells.Item(107, 7).Value = 'USD'
$ws.Cells.Item(107, 8).Formula = "=SUM(F107*I107*-1)"
$ws.Cells.Item(107, 9).Value = 1.1304
$ws.Cells.Item(107, 9).Style = "Normal"
$ws.Cells.Item(107, 11).Value = 'Barcap'
$ws.Cells.Item(107, 12).Value = 42509
$ws.Cells.Item(107, 13).Value = 'Louella Spiker'
$ws.Cells.Item(107, 14).Value = 'Approved'

$ws.Cells.Item(108, 1).Value = 107
$ws.Cells.Item(108, 2).Value = 42593
$ws.Cells.Item(108, 3).Formula = "=SUM(B108+1)"
$ws.Cells.Item(108, 4).Value = 'Spot'
$ws.Cells.Item(108, 5).Value = 'GBP'
$ws.Cells.Item(108, 6).Value = 500000
$ws.Cells.Item(108, 7).Value = 'USD'
$ws.Cells.Item(108, 8).Formula = "=SUM(F108*I108*-1)"
$ws.Cells.Item(108, 9).Value = 1.4352
$ws.Cells.Item(108, 9).Style = "Normal"
$ws.Cells.Item(108, 11).Value = 'Royal Bank of Scotland'
$ws.Cells.Item(108, 12).Value = 42509
$ws.Cells.Item(108, 13).Value = 'Magen Willison'
$ws.Cells.Item(108, 14).Value = 'Approved'

$ws.Cells.Item(109, 1).Value = 108
$ws.Cells.Item(109, 2).Value = 42634
$ws.Cells.Item(109, 3).Formula = "=SUM(B109+1)"
$ws.Cells.Item(109, 4).Value = 'Swap'
$ws.Cells.Item(109, 5).Value = 'EUR'
$ws.Cells.Item(109, 6).Value = 500000
$ws.Cells.Item(109, 7).Value = 'USD'
$ws.Cells.Item(109, 8).Formula = "=SUM(F109*I109*-1)"
$ws.Cells.Item(109, 9).Value = 1.1304
$ws.Cells.Item(109, 9).Style = "Normal"
$ws.Cells.Item(109, 11).Value = 'JP Morgan Chase'
$ws.Cells.Item(109, 12).Value = 42509
$ws.Cells.Item(109, 13).Value = 'Annemarie Rybicki'
$ws.Cells.Item(109, 14).Value = 'Approved'

$ws.Cells.Item(110, 1).Value = 109
$ws.Cells.Item(110, 2).Value = 42654
$ws.Cells.Item(110, 3).Formula = "=SUM(B110+1)"
$ws.Cells.Item(110, 4).Value = 'Swap'
$ws.Cells.Item(110, 5).Value = 'EUR'
$ws.Cells.Item(110, 6).Value = 1000000
$ws.Cells.Item(110, 7).Value = 'USD'
$ws.Cells.Item(110, 8).Formula = "=SUM(F110*I110*-1)"
$ws.Cells.Item(110, 9).Value = 1.1304
$ws.Cells.Item(110, 9).Style = "Normal"
$ws.Cells.Item(110, 11).Value = 'MUFJ'
$ws.Cells.Item(110, 12).Value = 42509
$ws.Cells.Item(110, 13).Value = 'Yael Rich'
$ws.Cells.Item(110, 14).Value = 'Approved'

$ws.Cells.Item(111, 1).Value = 110
$ws.Cells.Item(111, 2).Value = 42662
$ws.Cells.Item(111, 3).Formula = "=SUM(B111+1)"
$ws.Cells.Item(111, 4).Value = 'Forward'
$ws.Cells.Item(111, 5).Value = 'EUR'
$ws.Cells.Item(111, 6).Value = -1000000
$ws.Cells.Item(111, 7).Value = 'GBP'
$ws.Cells.Item(111, 8).Formula = "=SUM(F111*I111*-1)"
$ws.Cells.Item(111, 9).Value = 0.7876
$ws.Cells.Item(111, 9).Style = "Normal"
$ws.Cells.Item(111, 11).Value = 'Lloyds TSB'
$ws.Cells.Item(111, 12).Value = 42509
$ws.Cells.Item(111, 13).Value = 'Sarai Pilgrim'
$ws.Cells.Item(111, 14).Value = 'Approved'

$ws.Cells.Item(112, 1).Value = 111
$ws.Cells.Item(112, 2).Value = 42674
$ws.Cells.Item(112, 3).Formula = "=SUM(B112+1)"
$ws.Cells.Item(112, 4).Value = 'Forward'
$ws.Cells.Item(112, 5).Value = 'EUR'
$ws.Cells.Item(112, 6).Value = -10000
$ws.Cells.Item(112, 7).Value = 'GBP'
$ws.Cells.Item(112, 8).Formula = "=SUM(F112*I112*-1)"
$ws.Cells.Item(112, 9).Value = 0.7876
$ws.Cells.Item(112, 9).Style = "Normal"
$ws.Cells.Item(112, 11).Value = 'Societe Generale'
$ws.Cells.Item(112, 12).Value = 42509
$ws.Cells.Item(112, 13).Value = 'Sheba Dowdy'
$ws.Cells.Item(112, 14).Value = 'Pending'

$ws.Cells.Item(113, 1).Value = 112
$ws.Cells.Item(113, 2).Value = 42679
$ws.Cells.Item(113, 3).Formula = "=SUM(B113+1)"
$ws.Cells.Item(113, 4).Value = 'Swap'
$ws.Cells.Item(113, 5).Value = 'EUR'
$ws.Cells.Item(113, 6).Value = 10000
$ws.Cells.Item(113, 7).Value = 'USD'
$ws.Cells.Item(113, 8).Formula = "=SUM(F113*I113*-1)"
$ws.Cells.Item(113, 9).Value = 1.1304
$ws.Cells.Item(113, 9).Style = "Normal"
$ws.Cells.Item(113, 11).Value = 'Bank of America Merrill Lynch'
$ws.Cells.Item(113, 12).Value = 42509
$ws.Cells.Item(113, 13).Value = 'Sanjuana Kimsey'
$ws.Cells.Item(113, 14).Value = 'Approved'

$ws.Cells.Item(114, 1).Value = 113
$ws.Cells.Item(114, 2).Value = 42688
$ws.Cells.Item(114, 3).Formula = "=SUM(B114+1)"
$ws.Cells.Item(114, 4).Value = 'Swap'
$ws.Cells.Item(114, 5).Value = 'EUR'
$ws.Cells.Item(114, 6).Value = 500000
$ws.Cells.Item(114, 7).Value = 'USD'
$ws.Cells.Item(114, 8).Formula = "=SUM(F114*I114*-1)"
$ws.Cells.Item(114, 9).Value = 1.2962
$ws.Cells.Item(114, 9).Style = "Normal"
$ws.Cells.Item(114, 11).Value = 'Deutsche Bank'
$ws.Cells.Item(114, 12).Value = 42509
$ws.Cells.Item(114, 13).Value = 'Shante Hey'
$ws.Cells.Item(114, 14).Value = 'Approved'

$ws.Cells.Item(115, 1).Value = 114
$ws.Cells.Item(115, 2).Value = 42703
$ws.Cells.Item(115, 3).Formula = "=SUM(B115+1)"
$ws.Cells.Item(115, 4).Value = 'Forward'
$ws.Cells.Item(115, 5).Value = 'EUR'
$ws.Cells.Item(115, 6).Value = -500000
$ws.Cells.Item(115, 7).Value = 'GBP'
$ws.Cells.Item(115, 8).Formula = "=SUM(F115*I115*-1)"
$ws.Cells.Item(115, 9).Value = 0.7876
$ws.Cells.Item(115, 9).Style = "Normal"
$ws.Cells.Item(115, 11).Value = 'Deutsche Bank'
$ws.Cells.Item(115, 12).Value = 42509
$ws.Cells.Item(115, 13).Value = 'Shante Hey'
$ws.Cells.Item(115, 14).Value = 'Pending'

$ws.Cells.Item(116, 1).Value = 115
$ws.Cells.Item(116, 2).Value = 42708
$ws.Cells.Item(116, 3).Formula = "=SUM(B116+1)"
$ws.Cells.Item(116, 4).Value = 'Spot'
$ws.Cells.Item(116, 5).Value = 'GBP'
$ws.Cells.Item(116, 6).Value = 500000
$ws.Cells.Item(116, 7).Value = 'USD'
$ws.Cells.Item(116, 8).Formula = "=SUM(F116*I116*-1)"
$ws.Cells.Item(116, 9).Value = 1.4352
$ws.Cells.Item(116, 9).Style = "Normal"
$ws.Cells.Item(116, 11).Value = 'JP Morgan Chase'
$ws.Cells.Item(116, 12).Value = 42509
$ws.Cells.Item(116, 13).Value = 'Granville Westfall'
$ws.Cells.Item(116, 14).Value = 'Approved'

$ws.Cells.Item(117, 1).Value = 116
$ws.Cells.Item(117, 2).Value = 42716
$ws.Cells.Item(117, 3).Formula = "=SUM(B117+1)"
$ws.Cells.Item(117, 4).Value = 'Forward'
$ws.Cells.Item(117, 5).Value = 'EUR'
$ws.Cells.Item(117, 6).Value = -1000000
$ws.Cells.Item(117, 7).Value = 'GBP'
$ws.Cells.Item(117, 8).Formula = "=SUM(F117*I117*-1)"
$ws.Cells.Item(117, 9).Value = 0.9865
$ws.Cells.Item(117, 9).Style = "Normal"
$ws.Cells.Item(117, 11).Value = 'MUFJ'
$ws.Cells.Item(117, 12).Value = 42509
$ws.Cells.Item(117, 13).Value = 'Yael Rich'
$ws.Cells.Item(117, 14).Value = 'Pending'

$ws.Cells.Item(118, 1).Value = 117
$ws.Cells.Item(118, 2).Value = 42723
$ws.Cells.Item(118, 3).Formula = "=SUM(B118+1)"
$ws.Cells.Item(118, 4).Value = 'Forward'
$ws.Cells.Item(118, 5).Value = 'EUR'
$ws.Cells.Item(118, 6).Value = -50000000
$ws.Cells.Item(118, 7).Value = 'JPY'
$ws.Cells.Item(118, 8).Formula = "=SUM(F118*I118*-1)"
$ws.Cells.Item(118, 9).Value = 123.9075
$ws.Cells.Item(118, 9).Style = "Normal"
$ws.Cells.Item(118, 11).Value = 'UBS'
$ws.Cells.Item(118, 12).Value = 42509
$ws.Cells.Item(118, 13).Value = 'Bradley Chumley'
$ws.Cells.Item(118, 14).Value = 'Approved'

$ws.Cells.Item(119, 1).Value = 118
$ws.Cells.Item(119, 2).Value = 42739
$ws.Cells.Item(119, 3).Formula = "=SUM(B119+1)"
$ws.Cells.Item(119, 4).Value = 'Forward'
$ws.Cells.Item(119, 5).Value = 'EUR'
$ws.Cells.Item(119, 6).Value = -10000
$ws.Cells.Item(119, 7).Value = 'GBP'
$ws.Cells.Item(119, 8).Formula = "=SUM(F119*I119*-1)"
$ws.Cells.Item(119, 9).Value = 0.7876
$ws.Cells.Item(119, 9).Style = "Normal"
$ws.Cells.Item(119, 11).Value = 'Societe Generale'
$ws.Cells.Item(119, 12).Value = 42509
$ws.Cells.Item(119, 13).Value = 'Loralee Stalker'
$ws.Cells.Item(119, 14).Value = 'Pending'

$ws.Cells.Item(120, 1).Value = 119
$ws.Cells.Item(120, 2).Value = 42754
$ws.Cells.Item(120, 3).Formula = "=SUM(B120+1)"
$ws.Cells.Item(120, 4).Value = 'Forward'
$ws.Cells.Item(120, 5).Value = 'EUR'
$ws.Cells.Item(120, 6).Value = -50000000
$ws.Cells.Item(120, 7).Value = 'JPY'
$ws.Cells.Item(120, 8).Formula = "=SUM(F120*I120*-1)"
$ws.Cells.Item(120, 9).Value = 121.5689
$ws.Cells.Item(120, 9).Style = "Normal"
$ws.Cells.Item(120, 11).Value = 'Credit Suisse'
$ws.Cells.Item(120, 12).Value = 42517
$ws.Cells.Item(120, 13).Value = 'Britany Saffell'
$ws.Cells.Item(120, 14).Value = 'Approved'

$ws.Cells.Item(121, 1).Value = 120
$ws.Cells.Item(121, 2).Value = 42777
$ws.Cells.Item(121, 3).Formula = "=SUM(B121+1)"
$ws.Cells.Item(121, 4).Value = 'Spot'
$ws.Cells.Item(121, 5).Value = 'USD'
$ws.Cells.Item(121, 6).Value = 50000000
$ws.Cells.Item(121, 7).Value = 'JPY'
$ws.Cells.Item(121, 8).Formula = "=SUM(F121*I121*-1)"
$ws.Cells.Item(121, 9).Value = 109.6225
$ws.Cells.Item(121, 9).Style = "Normal"
$ws.Cells.Item(121, 11).Value = 'UBS'
$ws.Cells.Item(121, 12).Value = 42517
$ws.Cells.Item(121, 13).Value = 'Bradley Chumley'
$ws.Cells.Item(121, 14).Value = 'Approved'

$ws.Cells.Item(122, 1).Value = 121
$ws.Cells.Item(122, 2).Value = 42811
$ws.Cells.Item(122, 3).Formula = "=SUM(B122+1)"
$ws.Cells.Item(122, 4).Value = 'Spot'
$ws.Cells.Item(122, 5).Value = 'USD'
$ws.Cells.Item(122, 6).Value = 50000000
$ws.Cells.Item(122, 7).Value = 'JPY'
$ws.Cells.Item(122, 8).Formula = "=SUM(F122*I122*-1)"
$ws.Cells.Item(122, 9).Value = 109.3256
$ws.Cells.Item(122, 9).Style = "Normal"
$ws.Cells.Item(122, 11).Value = 'UBS'
$ws.Cells.Item(122, 12).Value = 42517
$ws.Cells.Item(122, 13).Value = 'Kelley Babb'
$ws.Cells.Item(122, 14).Value = 'Approved'

$ws.Cells.Item(123, 1).Value = 122
$ws.Cells.Item(123, 2).Value = 42844
$ws.Cells.Item(123, 3).Formula = "=SUM(B123+1)"
$ws.Cells.Item(123, 4).Value = 'Spot'
$ws.Cells.Item(123, 5).Value = 'USD'
$ws.Cells.Item(123, 6).Value = 10000
$ws.Cells.Item(123, 7).Value = 'JPY'
$ws.Cells.Item(123, 8).Formula = "=SUM(F123*I123*-1)"
$ws.Cells.Item(123, 9).Value = 107.6548
$ws.Cells.Item(123, 9).Style = "Normal"
$ws.Cells.Item(123, 11).Value = 'Bank of America Merrill Lynch'
$ws.Cells.Item(123, 12).Value = 42517
$ws.Cells.Item(123, 13).Value = 'Sanjuana Kimsey'
$ws.Cells.Item(123, 14).Value = 'Approved'

# --- Selection / view state ---
$ws.Range("B7").Select()
